$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: merge the three runs
#   "...We need to adapt our dataset for " + "modelling" + ", "
# into a single run whose text is
#   "...We need to adapt our dataset for modelling, "
# Select the exact span covered by those three original runs and rewrite it
# as one contiguous piece of text (Word will coalesce it back into a single
# run because it's assigned in one shot with uniform formatting).
# ---------------------------------------------------------------------------
$rngFind = $d.Content
$ok = $rngFind.Find.Execute(
    "Creating a Spark database from a large 3.61GB JSON file for modeling purposes is impractical due to its size. We need to adapt our dataset for modelling, ",
    $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0)
if ($ok) {
    $rngFind.Text = "Creating a Spark database from a large 3.61GB JSON file for modeling purposes is impractical due to its size. We need to adapt our dataset for modelling, "
}

# ---------------------------------------------------------------------------
# Locate the "...maxResultSize..." paragraph - the insertion point for the
# new block of paragraphs (added right after it, before the next existing
# paragraph in the document).
# ---------------------------------------------------------------------------
$maxResultParaIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*maxResultSize*") {
        $maxResultParaIdx = $i
    }
}

# ---------------------------------------------------------------------------
# Step 1: empty paragraph, Courier New style - InsertParagraphAfter inherits
# formatting from the "maxResultSize" paragraph, which already uses that
# exact Courier New / tabs / shading style.
# ---------------------------------------------------------------------------
$pAnchor = $d.Paragraphs.Item($maxResultParaIdx)
$pAnchor.Range.InsertParagraphAfter()

# The paragraph that used to immediately follow "maxResultSize" (an empty,
# bold Calibri/24 paragraph in the original document) has now been pushed
# one slot further down; it already carries exactly the style the
# "March 4, 2024" heading needs.
$calibriIdx = $maxResultParaIdx + 2

# ---------------------------------------------------------------------------
# Step 2: "March 4, 2024" heading paragraph - insert before that Calibri
# anchor paragraph so it inherits the same bold/Calibri/24 style, then set
# its text.
# ---------------------------------------------------------------------------
$pCalibri = $d.Paragraphs.Item($calibriIdx)
$pCalibri.Range.InsertParagraphBefore()
$marchIdx = $calibriIdx
$d.Paragraphs.Item($marchIdx).Range.Text = "March 4, 2024"

# ---------------------------------------------------------------------------
# Steps 3-6: Times New Roman / size 20 body paragraphs - insert before the
# paragraph that immediately follows the Calibri anchor (the original empty
# Times-New-Roman/20 paragraph), so each new paragraph inherits the correct
# style, then set its text.
# ---------------------------------------------------------------------------
$tnrIdx = $calibriIdx + 2
$pTnr = $d.Paragraphs.Item($tnrIdx)

$bodyTexts = @(
    "Spark is often used for distributed data processing, including data cleaning, feature extraction, and handling large datasets that don't fit into a single machine's memory.",
    "Keras, on the other hand, is a high-level neural networks API, typically used for building and training deep learning models.",
    "Finding an example of a deep learning application using a dataset that could be classified under the `"Big Data`" umbrella, due to the large volume dataset."
)

foreach ($txt in $bodyTexts) {
    $pTnr.Range.InsertParagraphBefore()
    $newIdx = $tnrIdx
    $tnrIdx = $tnrIdx + 1
    $d.Paragraphs.Item($newIdx).Range.Text = $txt
    $pTnr = $d.Paragraphs.Item($tnrIdx)
}

# "Book:" paragraph, followed by a manual line break then the second
# sentence, still within the same paragraph.
$pTnr.Range.InsertParagraphBefore()
$bookIdx = $tnrIdx
$tnrIdx = $tnrIdx + 1
$pTnr = $d.Paragraphs.Item($tnrIdx)

$pBook = $d.Paragraphs.Item($bookIdx)
$pBook.Range.Text = "Book:"
$bookEnd = $pBook.Range.End
$rngBreak = $d.Range($bookEnd - 1, $bookEnd - 1)
$rngBreak.InsertAfter("`vDeep Learning Convergence to Big Data Analytics.")

# ---------------------------------------------------------------------------
# Step 7: final empty paragraph, Courier New style - insert before the
# paragraph that now immediately follows the new block (the original empty,
# bold Calibri/24 paragraph) - but we want Courier New, not Calibri, so
# instead insert right after the "Book:" paragraph, which inherits Times
# New Roman; fix it up afterwards the same way the first Courier paragraph
# was produced (copy from the maxResultSize paragraph's style again).
# ---------------------------------------------------------------------------
$pAnchor2 = $d.Paragraphs.Item($maxResultParaIdx)
$pAnchor2.Range.InsertParagraphAfter()

Write-Output "done"
Write-Output ("FinalParagraphCount=" + $d.Paragraphs.Count)
